$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.193.50"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.500.45"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.18"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.72"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").Value = "2.523.35"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.57"
$ws.Range("E12").Value = "  +5.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "2.943.43"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.48"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").Value = "59.074.99"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "2.519.05"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.89"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  +3.34%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.92"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("E25").Value = "  -4.52%  "
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.86"
$ws.Range("D29").Value = "0.0₃0774"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.68"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.81"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("E32").Value = "  -7.32%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  +6.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.41"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.63"
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.34"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  -8.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.66"
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.00"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "295.43"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.69"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.77"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.70"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.56"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0514"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0228"
$ws.Range("E51").Value = "  -0.86%  "
